# Update "想去人数" (want-to-go count) figures scraped on a later run.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1650
$ws1.Range("F3").Value = 9056
$ws1.Range("F4").Value = 110
$ws1.Range("F5").Value = 503
$ws1.Range("F6").Value = 694
$ws1.Range("F7").Value = 777
$ws1.Range("F9").Value = 51
$ws1.Range("F10").Value = 84
$ws1.Range("F11").Value = 5666
$ws1.Range("F15").Value = 4351
$ws1.Range("F17").Value = 159
$ws1.Range("F19").Value = 16
$ws1.Range("F21").Value = 17
$ws1.Range("F24").Value = 2692
$ws1.Range("F25").Value = 122

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 8

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1650
$ws4.Range("F3").Value = 9056
$ws4.Range("F4").Value = 110
$ws4.Range("F5").Value = 8
$ws4.Range("F6").Value = 503
$ws4.Range("F7").Value = 694
$ws4.Range("F8").Value = 777
$ws4.Range("F10").Value = 51
$ws4.Range("F11").Value = 84
$ws4.Range("F12").Value = 5666
$ws4.Range("F16").Value = 4351
$ws4.Range("F18").Value = 159
$ws4.Range("F20").Value = 16
$ws4.Range("F22").Value = 17
$ws4.Range("F25").Value = 2692
$ws4.Range("F27").Value = 122
